$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "62.127.27"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "3.422.36"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'407.66"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").Value = "'133.29"
$ws.Range("E6").Value = "  +2.77%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.678"
$ws.Range("E9").Value = "  -2.16%  "
$ws.Range("D10").Value = "'0.122"
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("D11").Value = "'42.22"
$ws.Range("E11").Value = "  -3.78%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").Value = "3.967.46"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'8.44"
$ws.Range("E14").Value = "  -3.71%  "
$ws.Range("D15").Value = "'19.91"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").Value = "3.421.43"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "62.161.37"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "'0.0000132"
$ws.Range("E20").Value = "  -4.69%  "
$ws.Range("D21").Value = "'3.21"
$ws.Range("E21").Value = "  -4.09%  "
$ws.Range("D22").Value = "'85.03"
$ws.Range("E22").Value = "  +3.57%  "
$ws.Range("D23").Value = "'315.17"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").Value = "'12.81"
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("E25").Value = "  -2.87%  "
$ws.Range("D26").Value = "'4.78"
$ws.Range("E26").Value = "  +9.17%  "
$ws.Range("D27").Value = "'29.76"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").Value = "'8.29"
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("D29").Value = "'7.70"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("E30").Value = "  +3.72%  "
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("E32").Value = "  -5.42%  "
$ws.Range("D33").Value = "'42.96"
$ws.Range("E33").Value = "  -5.10%  "
$ws.Range("B34").Value = "Dai"
$ws.Range("C34").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("B35").Value = "Cosmos"
$ws.Range("C35").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D35").Value = "'11.40"
$ws.Range("E35").Value = "  -4.89%  "
$ws.Range("D36").Value = "'0.0485"
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("D37").Value = "'52.15"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("E39").Value = "  -4.45%  "
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "'137.87"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "'3.99"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "'16.85"
$ws.Range("E46").Value = "  -5.89%  "
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").Value = "'21.37"
$ws.Range("E48").Value = "  -5.80%  "
$ws.Range("D49").Value = "2.131.38"
$ws.Range("E49").Value = "  -5.24%  "
$ws.Range("D50").Value = "'2.29"
$ws.Range("E50").Value = "  -4.36%  "
$ws.Range("D51").Value = "'1.89"
$ws.Range("E51").Value = "  +1.36%  "
